# add the NA's under duplicate_image_filename
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
